$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the September row label to reflect the new "through" date
$ws.Range("A10").Value = "September (through 09-02)"

# Update September row values (row 10) for 2017, 2018, 2020, 2021
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 4
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 8

# Update Total row values (row 11) for 2017, 2018, 2020, 2021
$ws.Range("D11").Value = 559
$ws.Range("E11").Value = 494
$ws.Range("G11").Value = 789
$ws.Range("H11").Value = 1078
